$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1806.8572
$ws.Range("I17").Value2 = 1200
$ws.Range("J17").Value2 = 1875.8182
$ws.Range("K17").Value2 = 3600
$ws.Range("L17").Value2 = 5627.4546
$ws.Range("M17").Value2 = -3432
$ws.Range("N17").Value2 = -5963.4546
$ws.Range("H40").Value2 = 2423.24
$ws.Range("J40").Value2 = 2725.4119
$ws.Range("L40").Value2 = 2725.4119
$ws.Range("N40").Value2 = -3075.4119
$ws.Range("H100").Value2 = 8202.4
$ws.Range("I100").Value2 = 3250
$ws.Range("J100").Value2 = 11504
$ws.Range("K100").Value2 = 3250
$ws.Range("L100").Value2 = 11504
$ws.Range("M100").Value2 = -2709
$ws.Range("N100").Value2 = -12586
$ws.Range("H112").Value2 = 1711.9
$ws.Range("J112").Value2 = 1923
$ws.Range("L112").Value2 = 5769
$ws.Range("N112").Value2 = -7985
$ws.Range("H116").Value2 = 12768.28
$ws.Range("I116").Value2 = 14291.235
$ws.Range("K116").Value2 = 14291.235
$ws.Range("M116").Value2 = -10849.235
$ws.Range("H128").Value2 = 0
$ws.Range("J128").Value2 = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value2 = 0
$ws.Range("H132").Value2 = 61958.23
$ws.Range("I132").Value2 = 66640.72
$ws.Range("K132").Value2 = 199922.16
$ws.Range("M132").Value2 = -197392.16
$ws.Range("H138").Value2 = 2328.6365
$ws.Range("I138").Value2 = 1831.12
$ws.Range("J138").Value2 = 3883.375
$ws.Range("K138").Value2 = 5493.36
$ws.Range("L138").Value2 = 11650.125
$ws.Range("M138").Value2 = -353.3599999999997
$ws.Range("N138").Value2 = -21930.125

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value2 = 153713.86
$ws.Range("I6").Value2 = 10332.333
$ws.Range("J6").Value2 = 261250
$ws.Range("K6").Value2 = 10332.333
$ws.Range("L6").Value2 = 261250
$ws.Range("M6").Value2 = -10159.333
$ws.Range("N6").Value2 = -261596
$ws.Range("H32").Value2 = 3663.15
$ws.Range("I32").Value2 = 2434.0356
$ws.Range("J32").Value2 = 10116
$ws.Range("K32").Value2 = 2434.0356
$ws.Range("L32").Value2 = 10116
$ws.Range("M32").Value2 = -2147.0356
$ws.Range("N32").Value2 = -10690
$ws.Range("H61").Value2 = 6702704
$ws.Range("I61").Value2 = 6702704
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 6702704
$ws.Range("L61").Value2 = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value2 = -6702492
$ws.Range("H63").Value2 = 5451.846
$ws.Range("I63").Value2 = 2331.6667
$ws.Range("J63").Value2 = 6387.9
$ws.Range("K63").Value2 = 2331.6667
$ws.Range("L63").Value2 = 6387.9
$ws.Range("M63").Value2 = -1645.6667
$ws.Range("N63").Value2 = -7759.9
$ws.Range("H66").Value2 = 5451.846
$ws.Range("I66").Value2 = 2331.6667
$ws.Range("J66").Value2 = 6387.9
$ws.Range("K66").Value2 = 11658.3335
$ws.Range("L66").Value2 = 31939.5
$ws.Range("M66").Value2 = -8226.333500000001
$ws.Range("N66").Value2 = -38803.5
$ws.Range("H88").Value2 = 1890
$ws.Range("I88").Value2 = 1820
$ws.Range("J88").Value2 = 1960
$ws.Range("K88").Value2 = 1820
$ws.Range("L88").Value2 = 1960
$ws.Range("M88").Value2 = -1414
$ws.Range("N88").Value2 = -2772
$ws.Range("H91").Value2 = 1890
$ws.Range("I91").Value2 = 1820
$ws.Range("J91").Value2 = 1960
$ws.Range("K91").Value2 = 1820
$ws.Range("L91").Value2 = 1960
$ws.Range("M91").Value2 = -416
$ws.Range("N91").Value2 = -4768
$ws.Range("H109").Value2 = 70000
$ws.Range("J109").Value2 = 70000
$ws.Range("L109").Value2 = 70000
$ws.Range("N109").Value2 = -72774
$ws.Range("H112").Value2 = 10175.5
$ws.Range("I112").Value2 = 10351
$ws.Range("J112").Value2 = 10000
$ws.Range("K112").Value2 = 10351
$ws.Range("L112").Value2 = 10000
$ws.Range("M112").Value2 = -8874
$ws.Range("N112").Value2 = -12954
$ws.Range("H132").Value2 = 1944323.4
$ws.Range("I132").Value2 = 2587775.8
$ws.Range("K132").Value2 = 7763327.399999999
$ws.Range("M132").Value2 = -7760797.399999999
$ws.Range("H136").Value2 = 6702704
$ws.Range("I136").Value2 = 6702704
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 20108112
$ws.Range("L136").Value2 = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value2 = -20105562

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 2194.3684
$ws.Range("I20").Value2 = 2292.1538
$ws.Range("J20").Value2 = 1982.5
$ws.Range("K20").Value2 = 2292.1538
$ws.Range("L20").Value2 = 1982.5
$ws.Range("M20").Value2 = -2045.1538
$ws.Range("N20").Value2 = -2476.5
$ws.Range("H80").Value2 = 5451.143
$ws.Range("J80").Value2 = 2242.5454
$ws.Range("L80").Value2 = 2242.5454
$ws.Range("N80").Value2 = -4238.5454
$ws.Range("H83").Value2 = 5451.143
$ws.Range("J83").Value2 = 2242.5454
$ws.Range("L83").Value2 = 11212.727
$ws.Range("N83").Value2 = -21196.727

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 6913.5557
$ws.Range("I31").Value2 = 1990.091
$ws.Range("J31").Value2 = 10298.4375
$ws.Range("K31").Value2 = 1990.091
$ws.Range("L31").Value2 = 10298.4375
$ws.Range("M31").Value2 = -1695.091
$ws.Range("N31").Value2 = -10888.4375
$ws.Range("H34").Value2 = 6913.5557
$ws.Range("I34").Value2 = 1990.091
$ws.Range("J34").Value2 = 10298.4375
$ws.Range("K34").Value2 = 1990.091
$ws.Range("L34").Value2 = 10298.4375
$ws.Range("M34").Value2 = -1788.091
$ws.Range("N34").Value2 = -10702.4375
$ws.Range("H122").Value2 = 3337.4211
$ws.Range("I122").Value2 = 1319.7273
$ws.Range("K122").Value2 = 3959.1819
$ws.Range("M122").Value2 = -1509.1819
$ws.Range("H134").Value2 = 1884.6111
$ws.Range("I134").Value2 = 1884.6111
$ws.Range("J134").Value2 = 0
$ws.Range("K134").Value2 = 5653.8333
$ws.Range("L134").Value2 = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value2 = -3118.8333
$ws.Range("H135").Value2 = 0
$ws.Range("J135").Value2 = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value2 = 0

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value2 = 431650.6
$ws.Range("I14").Value2 = 431650.6
$ws.Range("K14").Value2 = 1294951.8
$ws.Range("M14").Value2 = -1294778.8
$ws.Range("H56").Value2 = 8447.083000000001
$ws.Range("I56").Value2 = 8447.083000000001
$ws.Range("K56").Value2 = 8447.083000000001
$ws.Range("M56").Value2 = -7917.083000000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 272106.38
$ws.Range("J11").Value2 = 558433.3
$ws.Range("L11").Value2 = 558433.3
$ws.Range("N11").Value2 = -558711.3
$ws.Range("H43").Value2 = 14024.889
$ws.Range("J43").Value2 = 22370
$ws.Range("L43").Value2 = 22370
$ws.Range("N43").Value2 = -22672
$ws.Range("H44").Value2 = 29649
$ws.Range("J44").Value2 = 29649
$ws.Range("L44").Value2 = 29649
$ws.Range("N44").Value2 = -30841
$ws.Range("H122").Value2 = 3570.5715
$ws.Range("I122").Value2 = 4802.6
$ws.Range("K122").Value2 = 14407.8
$ws.Range("M122").Value2 = -11957.8
$ws.Range("H132").Value2 = 3014321.5
$ws.Range("I132").Value2 = 4017430.8
$ws.Range("J132").Value2 = 4994
$ws.Range("K132").Value2 = 12052292.4
$ws.Range("L132").Value2 = 14982
$ws.Range("M132").Value2 = -12049762.4
$ws.Range("N132").Value2 = -20042
$ws.Range("H136").Value2 = 17144.688
$ws.Range("J136").Value2 = 17466
$ws.Range("L136").Value2 = 52398
$ws.Range("N136").Value2 = -57498

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 10633.6
$ws.Range("I16").Value2 = 10633.6
$ws.Range("K16").Value2 = 10633.6
$ws.Range("M16").Value2 = -10463.6
$ws.Range("H93").Value2 = 2798.4
$ws.Range("J93").Value2 = 2373
$ws.Range("L93").Value2 = 2373
$ws.Range("N93").Value2 = -4869
$ws.Range("H100").Value2 = 7877.579
$ws.Range("I100").Value2 = 2711.0667
$ws.Range("J100").Value2 = 27252
$ws.Range("K100").Value2 = 2711.0667
$ws.Range("L100").Value2 = 27252
$ws.Range("M100").Value2 = -2170.0667
$ws.Range("N100").Value2 = -28334
$ws.Range("H110").Value2 = 63321.875
$ws.Range("J110").Value2 = 63321.875
$ws.Range("L110").Value2 = 63321.875
$ws.Range("N110").Value2 = -71501.875
$ws.Range("H132").Value2 = 846266.75
$ws.Range("I132").Value2 = 1082865.2
$ws.Range("K132").Value2 = 3248595.6
$ws.Range("M132").Value2 = -3246065.6
$ws.Range("H136").Value2 = 4080.9375
$ws.Range("I136").Value2 = 3553.0334
$ws.Range("J136").Value2 = 11999.5
$ws.Range("K136").Value2 = 10659.1002
$ws.Range("L136").Value2 = 35998.5
$ws.Range("M136").Value2 = -8109.100199999999
$ws.Range("N136").Value2 = -41098.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value2 = 11000
$ws.Range("I52").Value2 = 11000
$ws.Range("K52").Value2 = 11000
$ws.Range("M52").Value2 = -10774
$ws.Range("H70").Value2 = 50084
$ws.Range("J70").Value2 = 50084
$ws.Range("L70").Value2 = 50084
$ws.Range("N70").Value2 = -50714
$ws.Range("H73").Value2 = 50084
$ws.Range("J73").Value2 = 50084
$ws.Range("L73").Value2 = 50084
$ws.Range("N73").Value2 = -52268
$ws.Range("H112").Value2 = 53387
$ws.Range("J112").Value2 = 53387
$ws.Range("L112").Value2 = 53387
$ws.Range("N112").Value2 = -56341
$ws.Range("H113").Value2 = 4572.8823
$ws.Range("I113").Value2 = 2415.7778
$ws.Range("K113").Value2 = 7247.3334
$ws.Range("M113").Value2 = -5077.3334
$ws.Range("H122").Value2 = 3288.2354
$ws.Range("I122").Value2 = 2744.7778
$ws.Range("J122").Value2 = 3899.625
$ws.Range("K122").Value2 = 8234.3334
$ws.Range("L122").Value2 = 11698.875
$ws.Range("M122").Value2 = -5784.3334
$ws.Range("N122").Value2 = -16598.875
$ws.Range("H132").Value2 = 4911814.5
$ws.Range("I132").Value2 = 6494377
$ws.Range("K132").Value2 = 19483131
$ws.Range("M132").Value2 = -19480601
$ws.Range("H136").Value2 = 9530976
$ws.Range("I136").Value2 = 11551501
$ws.Range("K136").Value2 = 34654503
$ws.Range("M136").Value2 = -34651953
